# Add a new "VSTAT License File" row to the Common sheet, right before the
# existing "SD-WAN Portal License File" row (new row 70), shifting every
# row from 70 downward by one and carrying its cell comment along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")

$insertRow = 70
$lastRow = 129

# --- 1. Snapshot the comments for rows insertRow..lastRow (column A) before
#        the row insert shifts cell contents (legacy comments are NOT moved
#        automatically by Rows.Insert in this host, so we handle them by
#        hand below).
$savedComments = @()
for ($r = $insertRow; $r -le $lastRow; $r++) {
    $cell = $ws.Range("A$r")
    $cmt = $cell.Comment
    if ($cmt -ne $null) {
        $savedComments += , @($r, $cmt.Text())
    }
}

# --- 2. Insert the new row; this shifts cell values/styles/merges/data
#        validations down by one automatically.
$ws.Rows.Item($insertRow).Insert()

# The freshly-inserted row doesn't reliably pick up the same cell style as
# its neighbours, so copy the formatting explicitly from the row right
# below it (the row that used to be the insertion row, now shifted down).
$ws.Range("A$($insertRow + 1):B$($insertRow + 1)").Copy() | Out-Null
$ws.Range("A$($insertRow):B$($insertRow)").PasteSpecial(-4122) | Out-Null

# --- 3. Re-home the saved comments one row down, processing from the
#        bottom up so we never overwrite a comment we still need to move.
for ($i = $savedComments.Count - 1; $i -ge 0; $i--) {
    $origRow = $savedComments[$i][0]
    $text = $savedComments[$i][1]
    $newRow = $origRow + 1

    $oldCell = $ws.Range("A$origRow")
    if ($oldCell.Comment -ne $null) {
        $oldCell.Comment.Delete()
    }
    $ws.Range("A$newRow").AddComment($text) | Out-Null
}

# --- 4. Populate the newly-inserted row.
$ws.Range("A$insertRow").Value = "VSTAT License File"
$ws.Range("A$insertRow").AddComment("Optional License File for Elasticsearch [default: ]") | Out-Null
